# Add "hydrogen combined cycle" as a new power plant type, and rename
# the existing "hydrogen" row to "hydrogen combustion turbine".

$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$ws = $wb.Worksheets.Item("SoCtMbCtbDP")

# Rename existing "hydrogen" entry (row 24) to be more specific.
$ws.Range("A24").Value = "hydrogen combustion turbine"

# Add the new "hydrogen combined cycle" row right after it.
$ws.Range("A25").Value = "hydrogen combined cycle"
$ws.Range("B25").Value = 1

# Move the active tab / selection back to the About sheet, matching the
# saved view state in the updated workbook.
$ws.Range("A29").Select()
$wsAbout.Activate()
